$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cell A1 mit wichtigen Informationen"
$ws.Range("A1:B1").NumberFormat = "0.00"
$ws.Range("A1:B1").Style = "Normal"
